# Update countries & provincias Spain
#
# The "Pais" sheet is a country ranking table (col A = country name,
# B..H = case statistics). This update:
#   1) Inserts three "new" countries into the ranking (Indonesia, Austria,
#      Oman), each landing one slot above where it used to sit, pushing the
#      country that used to occupy that slot (and the one below it, for the
#      Indonesia block) down by one row with their *previous* figures.
#   2) Refreshes a few rows' statistics in place (no reordering).
#   3) Bumps the "last updated" timestamp from 10:35 to 11:05.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 31: Banglades -- figures refreshed in place -----------------------
$ws.Range("B31").Value = 28511
$ws.Range("C31").Value = 1773
$ws.Range("D31").Value = 5602
$ws.Range("E31").Value = 22501
$ws.Range("G31").Value = 22
$ws.Range("H31").Value = 408

# --- Rows 34-36: Indonesia enters the table ahead of Polonia / Ucrania -----
$ws.Range("A34").Value = "Indonesia"
$ws.Range("B34").Value = 20162
$ws.Range("C34").Value = 973
$ws.Range("D34").Value = 4838
$ws.Range("E34").Value = 14046
$ws.Range("G34").Value = 36
$ws.Range("H34").Value = 1278

$ws.Range("A35").Value = "Polonia"
$ws.Range("B35").Value = 19983
$ws.Range("C35").Value = 244
$ws.Range("D35").Value = 8452
$ws.Range("E35").Value = 10566
$ws.Range("G35").Value = 3
$ws.Range("H35").Value = 965

$ws.Range("A36").Value = "Ucrania"
$ws.Range("B36").Value = 19706
$ws.Range("C36").Value = 476
$ws.Range("D36").Value = 6227
$ws.Range("E36").Value = 12900
$ws.Range("G36").Value = 15
$ws.Range("H36").Value = 579

# --- Rows 42-43: Austria enters the table ahead of Japon --------------------
$ws.Range("A42").Value = "Austria"
$ws.Range("B42").Value = 16404
$ws.Range("C42").Value = 51
$ws.Range("D42").Value = 14951
$ws.Range("E42").Value = 820
$ws.Range("H42").Value = 633

$ws.Range("A43").Value = "Japon"
$ws.Range("B43").Value = 16385
$ws.Range("D43").Value = 12286
$ws.Range("E43").Value = 3328
$ws.Range("H43").Value = 771

# --- Row 60: Malasia -- figures refreshed in place ---------------------------
$ws.Range("B60").Value = 7059
$ws.Range("C60").Value = 50
$ws.Range("D60").Value = 5796
$ws.Range("E60").Value = 1149

# --- Rows 64-65: Oman enters the table ahead of Ghana ------------------------
$ws.Range("A64").Value = "Oman"
$ws.Range("B64").Value = 6370
$ws.Range("C64").Value = 327
$ws.Range("D64").Value = 1821
$ws.Range("E64").Value = 4519
$ws.Range("H64").Value = 30

$ws.Range("A65").Value = "Ghana"
$ws.Range("B65").Value = 6269
$ws.Range("D65").Value = 1898
$ws.Range("E65").Value = 4340
$ws.Range("H65").Value = 31

# --- Timestamp footer ---------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 11:05"
